$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (was exported with the wrong/raw XPC name)
$ws.Name = "Gamma1F"

# Append a new averaged-intensity row (row 16) for the Gaussian-quadrature
# scheme "HexGrid-60degTilt5degRes" (reuses the shared string already used
# by row 15 / index 13).
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.786037861580783
$ws.Range("D16").Value = 0.173322833693635
$ws.Range("E16").Value = 1.081775024705596
$ws.Range("F16").Value = 1.786037861580783
$ws.Range("G16").Value = 0.563061020639455
$ws.Range("H16").Value = 1.308073162744521
$ws.Range("I16").Value = 1.164614654127239
$ws.Range("J16").Value = 0.173322833693635
$ws.Range("K16").Value = 0.6275489291996155
$ws.Range("L16").Value = 1.206793395390199
$ws.Range("M16").Value = 1.012814092915205

# Recomputed value for row 13 (G13) shifted in the last digit after the
# Gaussian-quadrature scheme was re-derived.
$ws.Range("G13").Value = 0.9876954109514365
